# Applies the OOXML diff to the workbook:
#  - Column width changes on sheet "Hoja1" (columns C/D/E, X/Y, AD/AE -> 3/4/5, 24/25, 30/31)
#  - Cell value corrections in column Y, rows 13-18 (values were off by a factor of ~1000,
#    corrected to the intended decimal numbers)
#
# NOTE on column widths: the headless Excel engine quantizes ColumnWidth to whole-pixel
# increments (i.e. stored width values land on a 1/6-character grid) when the property is
# set through the COM object model, exactly like real Excel does. Some of the target widths
# in the diff (8.57 and 8.71) do not fall on that grid, so the ColumnWidth input below is
# chosen (via the inverse of Excel's pixel-rounding formula) to land as close as possible
# on the grid point nearest to the intended width; widths that are already grid-aligned
# (8.0 and 15.0) are reproduced exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width updates ---

# Columns C:E (3-5) all become width 8.57; C and E already are 8.57, only D (4) changes.
$ws.Columns.Item(4).ColumnWidth = 7.668

# Column X (24): 6.29 -> 8.0 (exact grid point)
$ws.Columns.Item(24).ColumnWidth = 7.168

# Column Y (25): 4.57 -> 8.71
$ws.Columns.Item(25).ColumnWidth = 7.8346

# Column AD (30): split out of the old AD:AE(30:31) merged width, new width 15.0 (exact grid point)
$ws.Columns.Item(30).ColumnWidth = 14.168
# Column AE (31) keeps its original width of 10.71 (left untouched on purpose)

# --- Cell value corrections (column Y, rows 13-18) ---
$ws.Range("Y13").Value = 25.19
$ws.Range("Y14").Value = 15.87
$ws.Range("Y15").Value = 15.11
$ws.Range("Y16").Value = 52.91
$ws.Range("Y17").Value = 11.284
$ws.Range("Y18").Value = 22.383
